$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4757.5557
$ws.Range("I64").Value = 4702.5713
$ws.Range("K64").Value = 4702.5713
$ws.Range("M64").Value = -4454.5713
$ws.Range("H67").Value = 4757.5557
$ws.Range("I67").Value = 4702.5713
$ws.Range("K67").Value = 4702.5713
$ws.Range("M67").Value = -3844.5713
$ws.Range("H70").Value = 2002.8125
$ws.Range("I70").Value = 1816.8
$ws.Range("J70").Value = 2312.8333
$ws.Range("K70").Value = 5450.4
$ws.Range("L70").Value = 6938.499899999999
$ws.Range("M70").Value = -5180.4
$ws.Range("N70").Value = -7478.499899999999
$ws.Range("H73").Value = 2002.8125
$ws.Range("I73").Value = 1816.8
$ws.Range("J73").Value = 2312.8333
$ws.Range("K73").Value = 5450.4
$ws.Range("L73").Value = 6938.499899999999
$ws.Range("M73").Value = -4514.4
$ws.Range("N73").Value = -8810.499899999999
$ws.Range("H74").Value = 4595.9
$ws.Range("I74").Value = 3633
$ws.Range("M74").Value = -2697
$ws.Range("H77").Value = 4595.9
$ws.Range("I77").Value = 3633
$ws.Range("K77").Value = 18165
$ws.Range("M77").Value = -13485
$ws.Range("H80").Value = 2658.8845
$ws.Range("I80").Value = 1554.7222
$ws.Range("J80").Value = 5143.25
$ws.Range("K80").Value = 4664.1666
$ws.Range("L80").Value = 15429.75
$ws.Range("M80").Value = -3666.1666
$ws.Range("N80").Value = -17425.75
$ws.Range("H83").Value = 2658.8845
$ws.Range("I83").Value = 1554.7222
$ws.Range("J83").Value = 5143.25
$ws.Range("K83").Value = 13992.4998
$ws.Range("L83").Value = 46289.25
$ws.Range("M83").Value = -9000.4998
$ws.Range("N83").Value = -56273.25
$ws.Range("H86").Value = 60609164
$ws.Range("I86").Value = 22223486
$ws.Range("J86").Value = 142864200
$ws.Range("K86").Value = 22223486
$ws.Range("L86").Value = 142864200
$ws.Range("M86").Value = -22222363
$ws.Range("N86").Value = -142866446
$ws.Range("H88").Value = 71572856
$ws.Range("I88").Value = 83501336
$ws.Range("K88").Value = 83501336
$ws.Range("M88").Value = -83500930
$ws.Range("H89").Value = 60609164
$ws.Range("I89").Value = 22223486
$ws.Range("J89").Value = 142864200
$ws.Range("K89").Value = 111117430
$ws.Range("L89").Value = 714321000
$ws.Range("M89").Value = -111111814
$ws.Range("N89").Value = -714332232
$ws.Range("H91").Value = 71572856
$ws.Range("I91").Value = 83501336
$ws.Range("K91").Value = 83501336
$ws.Range("M91").Value = -83499932
$ws.Range("H106").Value = 372924
$ws.Range("I106").Value = 696595.5
$ws.Range("K106").Value = 696595.5
$ws.Range("M106").Value = -695964.5
$ws.Range("H109").Value = 74226.5
$ws.Range("J109").Value = 74226.5
$ws.Range("L109").Value = 74226.5
$ws.Range("N109").Value = -77000.5
$ws.Range("H123").Value = 81730
$ws.Range("J123").Value = 82018.336
$ws.Range("L123").Value = 82018.336
$ws.Range("N123").Value = -91818.336
$ws.Range("H132").Value = 1808.3478
$ws.Range("I132").Value = 1510.3611
$ws.Range("J132").Value = 2881.1
$ws.Range("K132").Value = 4531.0833
$ws.Range("L132").Value = 8643.299999999999
$ws.Range("M132").Value = -2001.0833
$ws.Range("N132").Value = -13703.3
$ws.Range("H133").Value = 92395.62
$ws.Range("J133").Value = 92395.62
$ws.Range("L133").Value = 92395.62
$ws.Range("N133").Value = -102515.62
$ws.Range("H134").Value = 57898.57
$ws.Range("J134").Value = 57898.57
$ws.Range("L134").Value = 57898.57
$ws.Range("N134").Value = -68038.57000000001
$ws.Range("H136").Value = 78844.57000000001
$ws.Range("J136").Value = 78844.57000000001
$ws.Range("L136").Value = 78844.57000000001
$ws.Range("N136").Value = -89044.57000000001
$ws.Range("H137").Value = 323841.03
$ws.Range("I137").Value = 1740.6
$ws.Range("K137").Value = 5221.799999999999
$ws.Range("M137").Value = -2671.799999999999
$ws.Range("H139").Value = 99988.57000000001
$ws.Range("J139").Value = 99988.57000000001
$ws.Range("L139").Value = 99988.57000000001
$ws.Range("N139").Value = -110268.57
$ws.Range("H140").Value = 80556.5
$ws.Range("J140").Value = 80778.86
$ws.Range("L140").Value = 80778.86
$ws.Range("N140").Value = -91138.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 871.3333
$ws.Range("I2").Value = 700.5
$ws.Range("J2").Value = 1119.8182
$ws.Range("K2").Value = 700.5
$ws.Range("L2").Value = 1119.8182
$ws.Range("M2").Value = -587.5
$ws.Range("N2").Value = -1345.8182
$ws.Range("H32").Value = 16301.808
$ws.Range("I32").Value = 6332.7646
$ws.Range("J32").Value = 35132.223
$ws.Range("K32").Value = 6332.7646
$ws.Range("L32").Value = 35132.223
$ws.Range("M32").Value = -6045.7646
$ws.Range("N32").Value = -35706.223
$ws.Range("H45").Value = 2642.9333
$ws.Range("I45").Value = 2501.25
$ws.Range("J45").Value = 2804.8572
$ws.Range("K45").Value = 2501.25
$ws.Range("L45").Value = 2804.8572
$ws.Range("M45").Value = -2124.25
$ws.Range("N45").Value = -3558.8572
$ws.Range("H74").Value = 2025
$ws.Range("I74").Value = 1435.9048
$ws.Range("K74").Value = 1435.9048
$ws.Range("M74").Value = -561.9048
$ws.Range("H77").Value = 2025
$ws.Range("I77").Value = 1435.9048
$ws.Range("K77").Value = 7179.524
$ws.Range("M77").Value = -2811.524
$ws.Range("H88").Value = 753.7778
$ws.Range("I88").Value = 483.33334
$ws.Range("J88").Value = 889
$ws.Range("K88").Value = 483.33334
$ws.Range("L88").Value = 889
$ws.Range("M88").Value = -77.33334000000002
$ws.Range("N88").Value = -1701
$ws.Range("H91").Value = 753.7778
$ws.Range("I91").Value = 483.33334
$ws.Range("J91").Value = 889
$ws.Range("K91").Value = 483.33334
$ws.Range("L91").Value = 889
$ws.Range("M91").Value = 920.66666
$ws.Range("N91").Value = -3697
$ws.Range("H110").Value = 868.25
$ws.Range("I110").Value = 866.13336
$ws.Range("K110").Value = 866.13336
$ws.Range("M110").Value = 1178.86664
$ws.Range("H116").Value = 871.3333
$ws.Range("I116").Value = 700.5
$ws.Range("J116").Value = 1119.8182
$ws.Range("K116").Value = 700.5
$ws.Range("L116").Value = 1119.8182
$ws.Range("M116").Value = 1593.5
$ws.Range("N116").Value = -5707.8182
$ws.Range("H121").Value = 52402.5
$ws.Range("J121").Value = 52402.5
$ws.Range("L121").Value = 52402.5
$ws.Range("N121").Value = -55896.5
$ws.Range("H132").Value = 1610.7826
$ws.Range("I132").Value = 1246.3667
$ws.Range("J132").Value = 2294.0625
$ws.Range("K132").Value = 3739.1001
$ws.Range("L132").Value = 6882.1875
$ws.Range("M132").Value = -1209.1001
$ws.Range("N132").Value = -11942.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 28082.223
$ws.Range("J2").Value = 28082.223
$ws.Range("L2").Value = 28082.223
$ws.Range("N2").Value = -28308.223
$ws.Range("H3").Value = 871.3333
$ws.Range("I3").Value = 700.5
$ws.Range("J3").Value = 1119.8182
$ws.Range("K3").Value = 700.5
$ws.Range("L3").Value = 1119.8182
$ws.Range("M3").Value = -586.5
$ws.Range("N3").Value = -1347.8182
$ws.Range("H5").Value = 671
$ws.Range("J5").Value = 634.6667
$ws.Range("L5").Value = 634.6667
$ws.Range("N5").Value = -860.6667
$ws.Range("H80").Value = 45795.816
$ws.Range("I80").Value = 166841.17
$ws.Range("J80").Value = 403.8125
$ws.Range("K80").Value = 166841.17
$ws.Range("L80").Value = 403.8125
$ws.Range("M80").Value = -165843.17
$ws.Range("N80").Value = -2399.8125
$ws.Range("H81").Value = 24153.5
$ws.Range("J81").Value = 24615
$ws.Range("L81").Value = 24615
$ws.Range("N81").Value = -26737
$ws.Range("H83").Value = 45795.816
$ws.Range("I83").Value = 166841.17
$ws.Range("J83").Value = 403.8125
$ws.Range("K83").Value = 834205.8500000001
$ws.Range("L83").Value = 2019.0625
$ws.Range("M83").Value = -829213.8500000001
$ws.Range("N83").Value = -12003.0625
$ws.Range("H84").Value = 24153.5
$ws.Range("J84").Value = 24615
$ws.Range("L84").Value = 73845
$ws.Range("N84").Value = -84453
$ws.Range("H99").Value = 1243184.2
$ws.Range("I99").Value = 65178.938
$ws.Range("K99").Value = 65178.938
$ws.Range("M99").Value = -63680.938
$ws.Range("H107").Value = 3019.9546
$ws.Range("I107").Value = 2710.8235
$ws.Range("K107").Value = 2710.8235
$ws.Range("M107").Value = -790.8235
$ws.Range("H114").Value = 82210.27
$ws.Range("J114").Value = 83669.3
$ws.Range("L114").Value = 83669.3
$ws.Range("N114").Value = -92347.3
$ws.Range("H116").Value = 73774
$ws.Range("J116").Value = 73774
$ws.Range("L116").Value = 73774
$ws.Range("N116").Value = -82952
$ws.Range("H118").Value = 49709.5
$ws.Range("J118").Value = 49061
$ws.Range("L118").Value = 49061
$ws.Range("N118").Value = -52375
$ws.Range("H122").Value = 99982
$ws.Range("J122").Value = 99982
$ws.Range("L122").Value = 99982
$ws.Range("N122").Value = -109782
$ws.Range("H140").Value = 100665
$ws.Range("J140").Value = 64361.816
$ws.Range("N140").Value = -74721.81599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3722.923
$ws.Range("I31").Value = 2334
$ws.Range("J31").Value = 4913.4287
$ws.Range("K31").Value = 2334
$ws.Range("L31").Value = 4913.4287
$ws.Range("M31").Value = -2039
$ws.Range("N31").Value = -5503.4287
$ws.Range("H34").Value = 3722.923
$ws.Range("I34").Value = 2334
$ws.Range("J34").Value = 4913.4287
$ws.Range("K34").Value = 2334
$ws.Range("L34").Value = 4913.4287
$ws.Range("M34").Value = -2132
$ws.Range("N34").Value = -5317.4287
$ws.Range("H105").Value = 5116.1665
$ws.Range("I105").Value = 2499
$ws.Range("J105").Value = 6424.75
$ws.Range("K105").Value = 2499
$ws.Range("L105").Value = 6424.75
$ws.Range("M105").Value = -752
$ws.Range("N105").Value = -9918.75
$ws.Range("H107").Value = 1028.4286
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 1733
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 1733
$ws.Range("M107").Value = 1420
$ws.Range("N107").Value = -5573
$ws.Range("H116").Value = 51463.3
$ws.Range("J116").Value = 51463.3
$ws.Range("L116").Value = 51463.3
$ws.Range("N116").Value = -60641.3
$ws.Range("H117").Value = 40668.43
$ws.Range("J117").Value = 40668.43
$ws.Range("L117").Value = 40668.43
$ws.Range("N117").Value = -49846.43
$ws.Range("H119").Value = 59418.75
$ws.Range("J119").Value = 59418.75
$ws.Range("L119").Value = 59418.75
$ws.Range("N119").Value = -69094.75
$ws.Range("H122").Value = 2333.2222
$ws.Range("I122").Value = 2065.5417
$ws.Range("K122").Value = 6196.625100000001
$ws.Range("M122").Value = -3746.625100000001
$ws.Range("H132").Value = 2118
$ws.Range("J132").Value = 3250
$ws.Range("L132").Value = 9750
$ws.Range("N132").Value = -14810
$ws.Range("N134").Value = -12228
$ws.Range("H134").Value = 2435.4285
$ws.Range("I134").Value = 2443.6667
$ws.Range("J134").Value = 2386
$ws.Range("K134").Value = 7331.000100000001
$ws.Range("L134").Value = 7158
$ws.Range("M134").Value = -4796.000100000001
$ws.Range("H138").Value = 99835.71000000001
$ws.Range("J138").Value = 99835.71000000001
$ws.Range("L138").Value = 99835.71000000001
$ws.Range("N138").Value = -110115.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 136.25
$ws.Range("J2").Value = 59.5
$ws.Range("L2").Value = 357
$ws.Range("N2").Value = -583
$ws.Range("H60").Value = 963.125
$ws.Range("I60").Value = 1458.3334
$ws.Range("J60").Value = 922.9729599999999
$ws.Range("K60").Value = 4375.0002
$ws.Range("L60").Value = 2768.91888
$ws.Range("M60").Value = -4124.0002
$ws.Range("N60").Value = -3270.91888
$ws.Range("H68").Value = 936.3333
$ws.Range("I68").Value = 766.3333
$ws.Range("K68").Value = 2298.9999
$ws.Range("M68").Value = -1487.9999
$ws.Range("H71").Value = 936.3333
$ws.Range("I71").Value = 766.3333
$ws.Range("K71").Value = 6896.9997
$ws.Range("M71").Value = -2840.9997
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30766
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32652
$ws.Range("H81").Value = 1590.6
$ws.Range("J81").Value = 2619.8
$ws.Range("L81").Value = 7859.400000000001
$ws.Range("N81").Value = -10105.4
$ws.Range("H84").Value = 1590.6
$ws.Range("J84").Value = 2619.8
$ws.Range("L84").Value = 23578.2
$ws.Range("N84").Value = -34810.2
$ws.Range("H107").Value = 510.66666
$ws.Range("I107").Value = 579.53845
$ws.Range("K107").Value = 1738.61535
$ws.Range("M107").Value = 181.38465
$ws.Range("H121").Value = 1526.56
$ws.Range("I121").Value = 1100.6666
$ws.Range("J121").Value = 1661.0526
$ws.Range("K121").Value = 3301.9998
$ws.Range("L121").Value = 4983.1578
$ws.Range("M121").Value = -1991.9998
$ws.Range("N121").Value = -7603.1578
$ws.Range("H132").Value = 3491.7273
$ws.Range("I132").Value = 2575.3333
$ws.Range("J132").Value = 3835.375
$ws.Range("K132").Value = 23177.9997
$ws.Range("L132").Value = 34518.375
$ws.Range("M132").Value = -20647.9997
$ws.Range("N132").Value = -39578.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M20").Value = -39760
$ws.Range("H20").Value = 50159
$ws.Range("I20").Value = 40005
$ws.Range("K20").Value = 40005
$ws.Range("H24").Value = 27814.666
$ws.Range("J24").Value = 27814.666
$ws.Range("L24").Value = 27814.666
$ws.Range("N24").Value = -28160.666
$ws.Range("H93").Value = 79099.5
$ws.Range("J93").Value = 79099.5
$ws.Range("L93").Value = 79099.5
$ws.Range("N93").Value = -82843.5
$ws.Range("H122").Value = 18572.691
$ws.Range("I122").Value = 30255.285
$ws.Range("J122").Value = 4943
$ws.Range("K122").Value = 90765.855
$ws.Range("L122").Value = 14829
$ws.Range("M122").Value = -88315.855
$ws.Range("N122").Value = -19729
$ws.Range("H132").Value = 4053.641
$ws.Range("I132").Value = 3767.5293
$ws.Range("K132").Value = 11302.5879
$ws.Range("M132").Value = -8772.5879
$ws.Range("H140").Value = 70992.5
$ws.Range("J140").Value = 74277.14
$ws.Range("N140").Value = -84637.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 18012
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H61").Value = 1100.6
$ws.Range("I61").Value = 626.625
$ws.Range("K61").Value = 626.625
$ws.Range("M61").Value = -424.625
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 48000
$ws.Range("L63").Value = 48000
$ws.Range("N63").Value = -49498
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 48000
$ws.Range("L66").Value = 144000
$ws.Range("N66").Value = -151488
$ws.Range("H93").Value = 1783.8667
$ws.Range("I93").Value = 1614.1818
$ws.Range("J93").Value = 2250.5
$ws.Range("K93").Value = 1614.1818
$ws.Range("L93").Value = 2250.5
$ws.Range("M93").Value = -366.1818000000001
$ws.Range("N93").Value = -4746.5
$ws.Range("H113").Value = 1100.6
$ws.Range("I113").Value = 626.625
$ws.Range("K113").Value = 626.625
$ws.Range("M113").Value = 1543.375
$ws.Range("H121").Value = 54544.285
$ws.Range("J121").Value = 54544.285
$ws.Range("L121").Value = 54544.285
$ws.Range("N121").Value = -58038.285
$ws.Range("H132").Value = 7793.234
$ws.Range("I132").Value = 10499.774
$ws.Range("J132").Value = 2549.3125
$ws.Range("K132").Value = 31499.322
$ws.Range("L132").Value = 7647.9375
$ws.Range("M132").Value = -28969.322
$ws.Range("N132").Value = -12707.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10138.2
$ws.Range("J62").Value = 13007.7
$ws.Range("L62").Value = 13007.7
$ws.Range("N62").Value = -14255.7
$ws.Range("H65").Value = 10138.2
$ws.Range("J65").Value = 13007.7
$ws.Range("L65").Value = 65038.5
$ws.Range("N65").Value = -71278.5
$ws.Range("H96").Value = 2744.4443
$ws.Range("I96").Value = 2600
$ws.Range("K96").Value = 2600
$ws.Range("M96").Value = -1227
$ws.Range("H100").Value = 3402342
$ws.Range("I100").Value = 7937230.5
$ws.Range("J100").Value = 1175.75
$ws.Range("K100").Value = 15874461
$ws.Range("L100").Value = 2351.5
$ws.Range("M100").Value = -15873920
$ws.Range("N100").Value = -3433.5
$ws.Range("H103").Value = 49799.668
$ws.Range("J103").Value = 49799.668
$ws.Range("L103").Value = 49799.668
$ws.Range("N103").Value = -52143.668
$ws.Range("H122").Value = 4738.0586
$ws.Range("I122").Value = 4404.273
$ws.Range("J122").Value = 5350
$ws.Range("K122").Value = 13212.819
$ws.Range("L122").Value = 16050
$ws.Range("M122").Value = -10762.819
$ws.Range("N122").Value = -20950
$ws.Range("H126").Value = 2106.7778
$ws.Range("I126").Value = 2039.5294
$ws.Range("K126").Value = 6118.5882
$ws.Range("M126").Value = -3648.5882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N43").ClearContents()
